$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = 100951327
$ws.Range("B3").Value = 99398
$ws.Range("L3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("S3").Value = 10
$ws.Range("Y3").Value = "'2022-05-19"
$ws.Range("Z3").Value = "09:00"
$ws.Range("AA3").Value = "'2022-05-19"
$ws.Range("AB3").Value = "10:00"
$ws.Range("AC3").Value = "En blomma utslagen, övriga i knopp."
$ws.Range("AF3").Value = "'"
$ws.Range("AI3").Value = "Lövlund"

# --- Row 8 ---
$ws.Range("A8").Value = 100951332
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("I8").Value = "'"
$ws.Range("J8").Value = "'"
$ws.Range("K8").Value = "blomknopp"
$ws.Range("N8").Value = "'"
$ws.Range("S8").Value = 10
$ws.Range("Y8").Value = "'2022-05-19"
$ws.Range("Z8").Value = "09:00"
$ws.Range("AA8").Value = "'2022-05-19"
$ws.Range("AB8").Value = "10:00"
$ws.Range("AC8").Value = ""
$ws.Range("AF8").Value = "'"
